$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2275
$ws.Range("I70").Value = 1750
$ws.Range("J70").Value = 2800
$ws.Range("K70").Value = 5250
$ws.Range("L70").Value = 8400
$ws.Range("M70").Value = -4980
$ws.Range("N70").Value = -8940
# Row 73
$ws.Range("H73").Value = 2275
$ws.Range("I73").Value = 1750
$ws.Range("J73").Value = 2800
$ws.Range("K73").Value = 5250
$ws.Range("L73").Value = 8400
$ws.Range("M73").Value = -4314
$ws.Range("N73").Value = -10272
# Row 76
$ws.Range("H76").Value = 5941.9165
$ws.Range("I76").Value = 5126.3335
$ws.Range("J76").Value = 6524.476
$ws.Range("K76").Value = 5126.3335
$ws.Range("L76").Value = 6524.476
$ws.Range("M76").Value = -4811.3335
$ws.Range("N76").Value = -7154.476
# Row 79
$ws.Range("H79").Value = 5941.9165
$ws.Range("I79").Value = 5126.3335
$ws.Range("J79").Value = 6524.476
$ws.Range("K79").Value = 5126.3335
$ws.Range("L79").Value = 6524.476
$ws.Range("M79").Value = -4034.3335
$ws.Range("N79").Value = -8708.475999999999
# Row 92
$ws.Range("H92").Value = 2850.1428
$ws.Range("I92").Value = 2980.2
$ws.Range("K92").Value = 2980.2
$ws.Range("M92").Value = -1732.2
# Row 94
$ws.Range("H94").Value = 3621.3572
$ws.Range("I94").Value = 1962.375
$ws.Range("J94").Value = 5833.3335
$ws.Range("K94").Value = 1962.375
$ws.Range("L94").Value = 5833.3335
$ws.Range("M94").Value = -1511.375
$ws.Range("N94").Value = -6735.3335
# Row 138
$ws.Range("H138").Value = 2005.6765
$ws.Range("I138").Value = 1363.3182
$ws.Range("J138").Value = 3183.3333
$ws.Range("K138").Value = 4089.9546
$ws.Range("L138").Value = 9549.999899999999
$ws.Range("M138").Value = 1050.0454
$ws.Range("N138").Value = -19829.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7614.357
$ws.Range("I32").Value = 7712.778
$ws.Range("J32").Value = 4957
$ws.Range("K32").Value = 7712.778
$ws.Range("L32").Value = 4957
$ws.Range("M32").Value = -7425.778
$ws.Range("N32").Value = -5531

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 43303.668
$ws.Range("J75").Value = 43604
$ws.Range("L75").Value = 43604
$ws.Range("N75").Value = -45476
# Row 78
$ws.Range("H78").Value = 43303.668
$ws.Range("J78").Value = 43604
$ws.Range("L78").Value = 130812
$ws.Range("N78").Value = -140172
# Row 99
$ws.Range("H99").Value = 1021.6
$ws.Range("I99").Value = 902
$ws.Range("K99").Value = 902
$ws.Range("M99").Value = 596
# Row 105
$ws.Range("H105").Value = 3940.4443
$ws.Range("I105").Value = 2877.5
$ws.Range("J105").Value = 4526.8965
$ws.Range("K105").Value = 2877.5
$ws.Range("L105").Value = 4526.8965
$ws.Range("M105").Value = -1130.5
$ws.Range("N105").Value = -8020.8965

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9527966
$ws.Range("I31").Value = 4433.7417
$ws.Range("J31").Value = 83335336
$ws.Range("K31").Value = 4433.7417
$ws.Range("L31").Value = 83335336
$ws.Range("M31").Value = -4138.7417
$ws.Range("N31").Value = -83335926
# Row 34
$ws.Range("H34").Value = 9527966
$ws.Range("I34").Value = 4433.7417
$ws.Range("J34").Value = 83335336
$ws.Range("K34").Value = 4433.7417
$ws.Range("L34").Value = 83335336
$ws.Range("M34").Value = -4231.7417
$ws.Range("N34").Value = -83335740
# Row 47
$ws.Range("H47").Value = 40067.668
$ws.Range("J47").Value = 40067.668
$ws.Range("L47").Value = 40067.668
$ws.Range("N47").Value = -41199.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 1391.7646
$ws.Range("J86").Value = 1878.75
$ws.Range("L86").Value = 5636.25
$ws.Range("N86").Value = -8008.25
# Row 89
$ws.Range("H89").Value = 1391.7646
$ws.Range("J89").Value = 1878.75
$ws.Range("L89").Value = 16908.75
$ws.Range("N89").Value = -28764.75
# Row 95
$ws.Range("H95").Value = 11986.667
$ws.Range("J95").Value = 11986.667
$ws.Range("L95").Value = 35960.001
$ws.Range("N95").Value = -40078.001
# Row 122
$ws.Range("H122").Value = 956.5833
$ws.Range("I122").Value = 819.625
$ws.Range("K122").Value = 7376.625
$ws.Range("M122").Value = -4926.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 15000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 15000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -16640
# Row 80
$ws.Range("H80").Value = 12348379
$ws.Range("I80").Value = 20836028
$ws.Range("J80").Value = 2709
$ws.Range("K80").Value = 20836028
$ws.Range("L80").Value = 2709
$ws.Range("M80").Value = -20835030
$ws.Range("N80").Value = -4705
# Row 83
$ws.Range("H83").Value = 12348379
$ws.Range("I83").Value = 20836028
$ws.Range("J83").Value = 2709
$ws.Range("K83").Value = 104180140
$ws.Range("L83").Value = 13545
$ws.Range("M83").Value = -104175148
$ws.Range("N83").Value = -23529
# Row 97
$ws.Range("H97").Value = 1034.2858
$ws.Range("I97").Value = 1038.8889
$ws.Range("J97").Value = 1006.6667
$ws.Range("K97").Value = 1038.8889
$ws.Range("L97").Value = 1006.6667
$ws.Range("M97").Value = -542.8888999999999
$ws.Range("N97").Value = -1998.6667
# Row 126
$ws.Range("H126").Value = 3781.1035
$ws.Range("I126").Value = 2335.2856
$ws.Range("J126").Value = 5130.533
$ws.Range("K126").Value = 7005.8568
$ws.Range("L126").Value = 15391.599
$ws.Range("M126").Value = -4535.8568
$ws.Range("N126").Value = -20331.599
# Row 134
$ws.Range("H134").Value = 39163
$ws.Range("J134").Value = 39163
$ws.Range("L134").Value = 117489
$ws.Range("N134").Value = -122559

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4733.971
$ws.Range("I7").Value = 5360.6924
$ws.Range("K7").Value = 5360.6924
$ws.Range("M7").Value = -5248.6924
# Row 40
$ws.Range("H40").Value = 3545.457
$ws.Range("I40").Value = 4805.8667
$ws.Range("J40").Value = 2600.15
$ws.Range("K40").Value = 4805.8667
$ws.Range("L40").Value = 2600.15
$ws.Range("M40").Value = -4669.8667
$ws.Range("N40").Value = -2872.15
# Row 100
$ws.Range("H100").Value = 2120
$ws.Range("I100").Value = 1700
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 1700
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -1159
$ws.Range("N100").Value = -3482
# Row 126
$ws.Range("H126").Value = 4733.971
$ws.Range("I126").Value = 5360.6924
$ws.Range("K126").Value = 16082.0772
$ws.Range("M126").Value = -13612.0772
# Row 136
$ws.Range("H136").Value = 50018772
$ws.Range("I136").Value = 100002536
$ws.Range("J136").Value = 35001
$ws.Range("K136").Value = 300007608
$ws.Range("L136").Value = 105003
$ws.Range("M136").Value = -300005058
$ws.Range("N136").Value = -110103

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 923.3333
$ws.Range("I81").Value = 1152.6666
$ws.Range("J81").Value = 694
$ws.Range("K81").Value = 2305.3332
$ws.Range("L81").Value = 1388
$ws.Range("M81").Value = -1244.3332
$ws.Range("N81").Value = -3510
# Row 84
$ws.Range("H84").Value = 923.3333
$ws.Range("I84").Value = 1152.6666
$ws.Range("J84").Value = 694
$ws.Range("K84").Value = 11526.666
$ws.Range("L84").Value = 6940
$ws.Range("M84").Value = -6222.666000000001
$ws.Range("N84").Value = -17548
# Row 96
$ws.Range("H96").Value = 2225.4285
$ws.Range("I96").Value = 2183.5
$ws.Range("J96").Value = 2263.5454
$ws.Range("K96").Value = 2183.5
$ws.Range("L96").Value = 2263.5454
$ws.Range("M96").Value = -810.5
$ws.Range("N96").Value = -5009.5454
# Row 109
$ws.Range("H109").Value = 29377
$ws.Range("J109").Value = 29377
$ws.Range("L109").Value = 29377
$ws.Range("N109").Value = -32151
# Row 122
$ws.Range("H122").Value = 1746.862
$ws.Range("I122").Value = 1549.3478
$ws.Range("J122").Value = 2504
$ws.Range("K122").Value = 4648.0434
$ws.Range("L122").Value = 7512
$ws.Range("M122").Value = -2198.0434
$ws.Range("N122").Value = -12412
